$wb = $excel.ActiveWorkbook

# Update the "Date" value on the Metadata sheet
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-07-21T11:52:46+00:00"

# Update the "System URI" value on the Include #0 sheet
$wsInc0 = $wb.Worksheets.Item("Include #0")
$wsInc0.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-G13-OrientationParticuliere"

# Update the "System URI" value on the Include #1 sheet
$wsInc1 = $wb.Worksheets.Item("Include #1")
$wsInc1.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R01-EnsembleSavoirFaire-CISIS"
